# Wordix - estructurasDatosWordix.xlsx
# "Agregue el ultimo arreglo de palabras ya jugadas por un jugador"
#
# Adds the new $palabrasJugadas data-structure block (an indexed array of
# already-played words: QUESO, PIANO, CEJAS) right before the existing
# $estadistJugador block, refreshes the page heading, and cleans up the
# leftover placeholder text ("A continuacion..." / "**** COMPLETAR ****").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Page heading: was the intro blurb, now "Colecciones utilizadas"
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Colecciones utilizadas"

# ------------------------------------------------------------------
# 2) Remove the leftover placeholder rows (12-13) that invited the
#    student to add more structures - no longer needed now that the
#    new structure has actually been added.
# ------------------------------------------------------------------
$ws.Range("A12").Clear()
$ws.Range("A1").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A13").ClearContents()

# ------------------------------------------------------------------
# 3) Make room for the new $palabrasJugadas block: insert 9 fresh rows
#    right before the old row 46, pushing the $estadistJugador block
#    (and its trailing documentation) down to rows 55-61.
# ------------------------------------------------------------------
$ws.Rows("45:53").Insert()

# ------------------------------------------------------------------
# 4) Fill the new $palabrasJugadas structure (rows 45-53)
# ------------------------------------------------------------------

# Title row, formatted like the other collection headers ($estadistJugador, etc.)
$ws.Range("A36").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A45").Value = "$" + "palabrasJugadas"

# Index header row (0,1,2) formatted like the other index rows (e.g. $coleccionVocales)
$ws.Range("B28:D28").Copy()
$ws.Range("B45:D45").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B45").Value = 0
$ws.Range("C45").Value = 1
$ws.Range("D45").Value = 2

# Values row: the actual words played
$ws.Range("B4:D4").Copy()
$ws.Range("B46:D46").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B46").Value = "QUESO"
$ws.Range("C46").Value = "PIANO"
$ws.Range("D46").Value = "CEJAS"

# "Informacion de la estructura:" sub-header
$ws.Range("B22").Copy()
$ws.Range("B48").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B48").Value = "Informaci" + [char]0x00F3 + "n de la estructura:"

# Documentation lines (Tipo / Tipos de datos / Para que se utiliza)
$ws.Range("B40:B42").Copy()
$ws.Range("B49:B51").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B49").Value = "Tipo: Arreglo indexado"
$ws.Range("B50").Value = "Tipos de datos: Almacena datos de tipo String"
$ws.Range("B51").Value = [char]0x00BF + "Para qu" + [char]0x00E9 + " se utiliza?: Guardar  las palabras jugadas por un jugador determinado"

# Row heights for 51 (matches "wrapped" taller rows seen elsewhere) and
# the two blank trailing rows of this block
$ws.Rows("51").RowHeight = 19.8
$ws.Rows("52").RowHeight = 19.8
$ws.Rows("53").RowHeight = 19.8

# ------------------------------------------------------------------
# 5) New note appended right after the (now shifted) $estadistJugador
#    documentation block, explaining that $palabrasJugadas is
#    generated on demand and not persisted.
# ------------------------------------------------------------------
$ws.Range("A36").Copy()
$ws.Range("B62").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B62").Value = "Se genera solo cuando se ejecuta la funci" + [char]0x00F3 + "n que lo requiere. No se guarda"

$excel.CutCopyMode = $false
